$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new date columns were appended to the mobility table (JM:JO),
# continuing the daily date series and adding the corresponding data
# points for each of the 5 rows.

# Copy the number format (date style) from the last existing date header
# cell (JL1) onto the three new header cells so they pick up the same
# cellXf (numFmtId 14, m/d/yyyy) instead of minting a brand-new style.
$ws.Range("JL1").Copy()
$ws.Range("JM1:JO1").PasteSpecial(-4122)

# New date header values (one day after the previous, consecutive dates).
$ws.Range("JM1").Value2 = 44109
$ws.Range("JN1").Value2 = 44110
$ws.Range("JO1").Value2 = 44111

# Row 2 (New York City transit)
$ws.Range("JM2").Value2 = 49.63
$ws.Range("JN2").Value2 = 50.57
$ws.Range("JO2").Value2 = 50.36

# Row 3 (San Francisco - Bay Area transit)
$ws.Range("JM3").Value2 = 34.83
$ws.Range("JN3").Value2 = 32.64
$ws.Range("JO3").Value2 = 38.88

# Row 4 (Baltimore City transit)
$ws.Range("JM4").Value2 = 65.04
$ws.Range("JN4").Value2 = 59.62
$ws.Range("JO4").Value2 = 60.75

# Row 5 (Phoenix transit)
$ws.Range("JM5").Value2 = 59.69
$ws.Range("JN5").Value2 = 59.14
$ws.Range("JO5").Value2 = 62.59

# Match the saved view state: scrolled right toward the new columns with
# JU24 as the active cell/selection.
$excel.ActiveWindow.ScrollColumn = 114
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("JU24").Select()
